$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2025
$ws.Range("I40").Value = 1840
$ws.Range("J40").Value = 2333.3333
$ws.Range("K40").Value = 1840
$ws.Range("L40").Value = 2333.3333
$ws.Range("M40").Value = -1665
$ws.Range("N40").Value = -2683.3333

$ws.Range("H62").Value = 4537.5
$ws.Range("I62").Value = 4050
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 4050
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -3426
$ws.Range("N62").Value = -7248

$ws.Range("H65").Value = 4537.5
$ws.Range("I65").Value = 4050
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 20250
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -17130
$ws.Range("N65").Value = -36240

$ws.Range("H112").Value = 2527.1482
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 2919.682
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 8759.045999999998
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -10975.046

$ws.Range("H116").Value = 2716.4614
$ws.Range("I116").Value = 1980.8
$ws.Range("J116").Value = 5168.6665
$ws.Range("K116").Value = 1980.8
$ws.Range("L116").Value = 5168.6665
$ws.Range("M116").Value = 1461.2
$ws.Range("N116").Value = -12052.6665

$ws.Range("H141").Value = 605.8049
$ws.Range("I141").Value = 545.95
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 1637.85
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 3542.15
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12524.889
$ws.Range("I2").Value = 1418.5
$ws.Range("J2").Value = 34737.668
$ws.Range("K2").Value = 1418.5
$ws.Range("L2").Value = 34737.668
$ws.Range("M2").Value = -1305.5
$ws.Range("N2").Value = -34963.668

$ws.Range("H61").Value = 1135.5625
$ws.Range("I61").Value = 986.975
$ws.Range("K61").Value = 986.975
$ws.Range("M61").Value = -774.975

$ws.Range("H74").Value = 1123.6904
$ws.Range("I74").Value = 792.2121
$ws.Range("J74").Value = 2339.111
$ws.Range("K74").Value = 792.2121
$ws.Range("L74").Value = 2339.111
$ws.Range("M74").Value = 81.78790000000004
$ws.Range("N74").Value = -4087.111

$ws.Range("H77").Value = 1123.6904
$ws.Range("I77").Value = 792.2121
$ws.Range("J77").Value = 2339.111
$ws.Range("K77").Value = 3961.0605
$ws.Range("L77").Value = 11695.555
$ws.Range("M77").Value = 406.9395000000004
$ws.Range("N77").Value = -20431.555

$ws.Range("H116").Value = 12524.889
$ws.Range("I116").Value = 1418.5
$ws.Range("J116").Value = 34737.668
$ws.Range("K116").Value = 1418.5
$ws.Range("L116").Value = 34737.668
$ws.Range("M116").Value = 875.5
$ws.Range("N116").Value = -39325.668

$ws.Range("H136").Value = 1135.5625
$ws.Range("I136").Value = 986.975
$ws.Range("K136").Value = 2960.925
$ws.Range("M136").Value = -410.9250000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12524.889
$ws.Range("I3").Value = 1418.5
$ws.Range("J3").Value = 34737.668
$ws.Range("K3").Value = 1418.5
$ws.Range("L3").Value = 34737.668
$ws.Range("M3").Value = -1304.5
$ws.Range("N3").Value = -34965.668

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 812.75
$ws.Range("I58").Value = 696.8889
$ws.Range("J58").Value = 1438.4
$ws.Range("K58").Value = 696.8889
$ws.Range("L58").Value = 1438.4
$ws.Range("M58").Value = -493.8889
$ws.Range("N58").Value = -1844.4

$ws.Range("H136").Value = 812.75
$ws.Range("I136").Value = 696.8889
$ws.Range("J136").Value = 1438.4
$ws.Range("K136").Value = 2090.6667
$ws.Range("L136").Value = 4315.200000000001
$ws.Range("M136").Value = 459.3332999999998
$ws.Range("N136").Value = -9415.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1842.9
$ws.Range("I5").Value = 2117.2666
$ws.Range("J5").Value = 1019.8
$ws.Range("K5").Value = 6351.7998
$ws.Range("L5").Value = 3059.4
$ws.Range("M5").Value = -6239.7998
$ws.Range("N5").Value = -3283.4

$ws.Range("H34").Value = 2188.889
$ws.Range("I34").Value = 1233.3334
$ws.Range("J34").Value = 2666.6667
$ws.Range("K34").Value = 3700.0002
$ws.Range("L34").Value = 8000.000100000001
$ws.Range("M34").Value = -3616.0002
$ws.Range("N34").Value = -8168.000100000001

$ws.Range("H39").Value = 3025.2856
$ws.Range("J39").Value = 2814
$ws.Range("L39").Value = 8442
$ws.Range("N39").Value = -9030

$ws.Range("H55").Value = 2428.4285
$ws.Range("I55").Value = 1002
$ws.Range("J55").Value = 2999
$ws.Range("K55").Value = 3006
$ws.Range("L55").Value = 8997
$ws.Range("M55").Value = -2829
$ws.Range("N55").Value = -9351

$ws.Range("H106").Value = 2513.0908
$ws.Range("J106").Value = 2513.0908
$ws.Range("L106").Value = 7539.2724
$ws.Range("N106").Value = -9431.2724

$ws.Range("H113").Value = 678.25806
$ws.Range("J113").Value = 699.1724
$ws.Range("L113").Value = 2097.5172
$ws.Range("N113").Value = -6437.5172

$ws.Range("H122").Value = 847.13794
$ws.Range("I122").Value = 717.4286
$ws.Range("J122").Value = 888.4091
$ws.Range("K122").Value = 6456.8574
$ws.Range("L122").Value = 7995.6819
$ws.Range("M122").Value = -4006.8574
$ws.Range("N122").Value = -12895.6819

$ws.Range("H134").Value = 3225.44
$ws.Range("I134").Value = 1815.6923
$ws.Range("J134").Value = 4752.6665
$ws.Range("K134").Value = 5447.0769
$ws.Range("L134").Value = 14257.9995
$ws.Range("M134").Value = -377.0769
$ws.Range("N134").Value = -24397.9995

$ws.Range("H135").Value = 1842.9
$ws.Range("I135").Value = 2117.2666
$ws.Range("J135").Value = 1019.8
$ws.Range("K135").Value = 19055.3994
$ws.Range("L135").Value = 9178.199999999999
$ws.Range("M135").Value = -16520.3994
$ws.Range("N135").Value = -14248.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2889
$ws.Range("I80").Value = 1771.5714
$ws.Range("K80").Value = 1771.5714
$ws.Range("M80").Value = -773.5714

$ws.Range("H83").Value = 2889
$ws.Range("I83").Value = 1771.5714
$ws.Range("K83").Value = 8857.857
$ws.Range("M83").Value = -3865.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4267.875
$ws.Range("I136").Value = 4779.5
$ws.Range("J136").Value = 2050.8333
$ws.Range("K136").Value = 14338.5
$ws.Range("L136").Value = 6152.499899999999
$ws.Range("M136").Value = -11788.5
$ws.Range("N136").Value = -11252.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 36006.332
$ws.Range("I28").Value = 20000
$ws.Range("J28").Value = 44009.5
$ws.Range("K28").Value = 20000
$ws.Range("L28").Value = 44009.5
$ws.Range("M28").Value = -19652
$ws.Range("N28").Value = -44705.5

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H132").Value = 2312.1396
$ws.Range("I132").Value = 2377.9714
$ws.Range("K132").Value = 7133.914199999999
$ws.Range("M132").Value = -4603.914199999999
